$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of usage log data (rows 30 and 31)
$newRows = @(
    @{ Row = 30; A = 45690.329386574071; B = 8; C = 6; D = 193; E = 380; F = 356; G = 384; H = 2664; I = 384; J = 1216; K = 119; L = 304; M = 30; N = 2949; O = 3791 },
    @{ Row = 31; A = 45690.951608796298; B = 8; C = 6; D = 193; E = 380; F = 356; G = 384; H = 2664; I = 384; J = 1216; K = 119; L = 304; M = 30; N = 2977; O = 3796 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2 = $r.A
    $ws.Cells.Item($row, 2).Value2 = $r.B
    $ws.Cells.Item($row, 3).Value2 = $r.C
    $ws.Cells.Item($row, 4).Value2 = $r.D
    $ws.Cells.Item($row, 5).Value2 = $r.E
    $ws.Cells.Item($row, 6).Value2 = $r.F
    $ws.Cells.Item($row, 7).Value2 = $r.G
    $ws.Cells.Item($row, 8).Value2 = $r.H
    $ws.Cells.Item($row, 9).Value2 = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
}
